$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44
$ws.Cells.Item(44, 4).Value = 44736
$ws.Cells.Item(44, 11).Value = "Fuyu"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 250
$ws.Cells.Item(44, 14).Value = 19000
$ws.Cells.Item(44, 15).Value = 20000
$ws.Cells.Item(44, 16).Value = 19400
$ws.Cells.Item(44, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(44, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(44, 19).Value = 1293
$ws.Cells.Item(44, 20).Value = 15

# Row 45
$ws.Cells.Item(45, 4).Value = 44299
$ws.Cells.Item(45, 11).Value = "Fuyu"
$ws.Cells.Item(45, 12).Value = "Primera"
$ws.Cells.Item(45, 13).Value = 45
$ws.Cells.Item(45, 14).Value = 22000
$ws.Cells.Item(45, 15).Value = 22000
$ws.Cells.Item(45, 16).Value = 22000
$ws.Cells.Item(45, 17).Value = "$/bandeja 15 kilos empedrada"
$ws.Cells.Item(45, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(45, 19).Value = 1467
$ws.Cells.Item(45, 20).Value = 15

# Row 46
$ws.Cells.Item(46, 4).Value = 44299
$ws.Cells.Item(46, 11).Value = "Mankaki"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 55
$ws.Cells.Item(46, 14).Value = 20000
$ws.Cells.Item(46, 15).Value = 20000
$ws.Cells.Item(46, 16).Value = 20000
$ws.Cells.Item(46, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(46, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(46, 19).Value = 1333
$ws.Cells.Item(46, 20).Value = 15

# Row 47
$ws.Cells.Item(47, 4).Value = 44706
$ws.Cells.Item(47, 11).Value = "Mankaki"
$ws.Cells.Item(47, 12).Value = "Primera"
$ws.Cells.Item(47, 13).Value = 45
$ws.Cells.Item(47, 14).Value = 18000
$ws.Cells.Item(47, 15).Value = 18000
$ws.Cells.Item(47, 16).Value = 18000
$ws.Cells.Item(47, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(47, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(47, 19).Value = 1000
$ws.Cells.Item(47, 20).Value = 18

# Row 48
$ws.Cells.Item(48, 4).Value = 44336
$ws.Cells.Item(48, 11).Value = "Fuyu"
$ws.Cells.Item(48, 12).Value = "Especial"
$ws.Cells.Item(48, 13).Value = 30
$ws.Cells.Item(48, 14).Value = 18000
$ws.Cells.Item(48, 15).Value = 18000
$ws.Cells.Item(48, 16).Value = 18000
$ws.Cells.Item(48, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(48, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(48, 19).Value = 1200
$ws.Cells.Item(48, 20).Value = 15

# Row 49
$ws.Cells.Item(49, 4).Value = 44336
$ws.Cells.Item(49, 11).Value = "Mankaki"
$ws.Cells.Item(49, 12).Value = "Especial"
$ws.Cells.Item(49, 13).Value = 50
$ws.Cells.Item(49, 14).Value = 18000
$ws.Cells.Item(49, 15).Value = 18000
$ws.Cells.Item(49, 16).Value = 18000
$ws.Cells.Item(49, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(49, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(49, 19).Value = 1200
$ws.Cells.Item(49, 20).Value = 15

# Row 50
$ws.Cells.Item(50, 4).Value = 44348
$ws.Cells.Item(50, 11).Value = "Hachiya"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 30
$ws.Cells.Item(50, 14).Value = 25000
$ws.Cells.Item(50, 15).Value = 25000
$ws.Cells.Item(50, 16).Value = 25000
$ws.Cells.Item(50, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(50, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(50, 19).Value = 1389
$ws.Cells.Item(50, 20).Value = 18

# Row 51
$ws.Cells.Item(51, 4).Value = 44348
$ws.Cells.Item(51, 11).Value = "Mankaki"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 40
$ws.Cells.Item(51, 14).Value = 25000
$ws.Cells.Item(51, 15).Value = 25000
$ws.Cells.Item(51, 16).Value = 25000
$ws.Cells.Item(51, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(51, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(51, 19).Value = 1389
$ws.Cells.Item(51, 20).Value = 18

# Row 52 (new row)
$ws.Cells.Item(52, 1).Value = 10
$ws.Cells.Item(52, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(52, 3).Value = "La Araucanía"
$ws.Cells.Item(52, 4).Value = 44307
$ws.Cells.Item(52, 4).NumberFormat = $ws.Cells.Item(51, 4).NumberFormat
$ws.Cells.Item(52, 5).Value = 9
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100107
$ws.Cells.Item(52, 8).Value = "Otros"
$ws.Cells.Item(52, 9).Value = 100107001
$ws.Cells.Item(52, 10).Value = "Caqui"
$ws.Cells.Item(52, 11).Value = "Fuyu"
$ws.Cells.Item(52, 12).Value = "Primera"
$ws.Cells.Item(52, 13).Value = 50
$ws.Cells.Item(52, 14).Value = 17000
$ws.Cells.Item(52, 15).Value = 17000
$ws.Cells.Item(52, 16).Value = 17000
$ws.Cells.Item(52, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(52, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(52, 19).Value = 1133
$ws.Cells.Item(52, 20).Value = 15
